# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2:F2, zh-cn!C2,
#   de-de!C2 all share this value).
# - Re-fit the Status/zh-cn/de-de columns that held that text now that the
#   text is shorter.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Columns shrink to fit the new (shorter) text.
$newWidth = 12.5
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
